{"js": "const replacements = [\n  [\"2024-11-24 Sunday\", \"2024-11-25 Monday\"],\n  [\"159\u00f77=22, 5\", \"270\u00f78=33, 6\"],\n  [\"109\u00f74=27, 1\", \"275\u00f73=91, 2\"],\n  [\"361\u00f72=180, 1\", \"539\u00f79=59, 8\"],\n  [\"534\u00f79=59, 3\", \"460\u00f78=57, 4\"],\n  [\"371\u00f75=74, 1\", \"103\u00f79=11, 4\"],\n  [\"853\u00f72=426, 1\", \"307\u00f74=76, 3\"],\n  [\"265\u00f74=66, 1\", \"905\u00f78=113, 1\"],\n  [\"246\u00f76=41, 0\", \"906\u00f79=100, 6\"],\n  [\"134\u00f75=26, 4\", \"488\u00f77=69, 5\"],\n  [\"912\u00f79=101, 3\", \"787\u00f75=157, 2\"],\n  [\"947\u00f73=315, 2\", \"205\u00f79=22, 7\"],\n  [\"331\u00f74=82, 3\", \"732\u00f77=104, 4\"],\n  [\"414\u00f79=46, 0\", \"149\u00f79=16, 5\"],\n  [\"364\u00f77=52, 0\", \"966\u00f76=161, 0\"],\n  [\"940\u00f72=470, 0\", \"214\u00f75=42, 4\"],\n  [\"816\u00f75=163, 1\", \"230\u00f78=28, 6\"],\n  [\"534\u00f72=267, 0\", \"755\u00f73=251, 2\"],\n  [\"873\u00f73=291, 0\", \"395\u00f75=79, 0\"],\n  [\"899\u00f78=112, 3\", \"299\u00f76=49, 5\"],\n  [\"400\u00f77=57, 1\", \"920\u00f76=153, 2\"],\n  [\"271\u00f74=67, 3\", \"122\u00f79=13, 5\"],\n  [\"547\u00f79=60, 7\", \"426\u00f74=106, 2\"],\n  [\"581\u00f78=72, 5\", \"173\u00f73=57, 2\"],\n  [\"287\u00f75=57, 2\", \"896\u00f79=99, 5\"],\n  [\"613\u00f73=204, 1\", \"810\u00f73=270, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('text');\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2024-11-24 Sunday\", \"2024-11-25 Monday\")\n    ,@(\"159\u00f77=22, 5\", \"270\u00f78=33, 6\")\n    ,@(\"109\u00f74=27, 1\", \"275\u00f73=91, 2\")\n    ,@(\"361\u00f72=180, 1\", \"539\u00f79=59, 8\")\n    ,@(\"534\u00f79=59, 3\", \"460\u00f78=57, 4\")\n    ,@(\"371\u00f75=74, 1\", \"103\u00f79=11, 4\")\n    ,@(\"853\u00f72=426, 1\", \"307\u00f74=76, 3\")\n    ,@(\"265\u00f74=66, 1\", \"905\u00f78=113, 1\")\n    ,@(\"246\u00f76=41, 0\", \"906\u00f79=100, 6\")\n    ,@(\"134\u00f75=26, 4\", \"488\u00f77=69, 5\")\n    ,@(\"912\u00f79=101, 3\", \"787\u00f75=157, 2\")\n    ,@(\"947\u00f73=315, 2\", \"205\u00f79=22, 7\")\n    ,@(\"331\u00f74=82, 3\", \"732\u00f77=104, 4\")\n    ,@(\"414\u00f79=46, 0\", \"149\u00f79=16, 5\")\n    ,@(\"364\u00f77=52, 0\", \"966\u00f76=161, 0\")\n    ,@(\"940\u00f72=470, 0\", \"214\u00f75=42, 4\")\n    ,@(\"816\u00f75=163, 1\", \"230\u00f78=28, 6\")\n    ,@(\"534\u00f72=267, 0\", \"755\u00f73=251, 2\")\n    ,@(\"873\u00f73=291, 0\", \"395\u00f75=79, 0\")\n    ,@(\"899\u00f78=112, 3\", \"299\u00f76=49, 5\")\n    ,@(\"400\u00f77=57, 1\", \"920\u00f76=153, 2\")\n    ,@(\"271\u00f74=67, 3\", \"122\u00f79=13, 5\")\n    ,@(\"547\u00f79=60, 7\", \"426\u00f74=106, 2\")\n    ,@(\"581\u00f78=72, 5\", \"173\u00f73=57, 2\")\n    ,@(\"287\u00f75=57, 2\", \"896\u00f79=99, 5\")\n    ,@(\"613\u00f73=204, 1\", \"810\u00f73=270, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
